# The PDFs linked from column A ("Paper") were moved from the "paper/"
# folder to "paper_old/". Point each paper hyperlink at its new location.
# Presenter names (col B) and the "Files"/presentation links (col C) are
# untouched - only the href (and matching visible path text) in col A changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-PaperLink($row) {
    $cell = $ws.Cells.Item($row, 1)
    $cur = $cell.Text
    $idx = $cur.IndexOf('paper/')
    if ($idx -ge 0) {
        $newVal = $cur.Substring(0, $idx) + 'paper_old/' + $cur.Substring($idx + 6)
        $cell.Value = $newVal
    }
}

# Update in the same order the rows were originally authored (Chen, Sims,
# Rambachan, Rigby, Hill, Wickramasuriya) rather than top-to-bottom sheet
# order, since that's the sequence the links were created/touched in.
Update-PaperLink 4
Update-PaperLink 3
Update-PaperLink 2
Update-PaperLink 5
Update-PaperLink 7
Update-PaperLink 6

# Move the active selection to match the saved cursor position.
$ws.Range("A12").Select() | Out-Null
